$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.477.21"
$ws.Range("E2").Value = "  +1.00%  "

$ws.Range("D3").Value = "2.214.05"
$ws.Range("E3").Value = "  -1.11%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.80"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.610"
$ws.Range("E6").Value = "  -1.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.23"
$ws.Range("E7").Value = "  +1.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  +0.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.22"
$ws.Range("E10").Value = "  -2.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0925"
$ws.Range("E11").Value = "  -2.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.93"
$ws.Range("E12").Value = "  -2.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.86"
$ws.Range("E13").Value = "  -0.60%  "

$ws.Range("E14").Value = "  -1.94%  "

$ws.Range("D15").Value = "2.543.36"
$ws.Range("E15").Value = "  -1.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.65"
$ws.Range("E16").Value = "  +2.24%  "

$ws.Range("D17").Value = "2.208.43"
$ws.Range("E17").Value = "  -2.21%  "

$ws.Range("E18").Value = "  -4.40%  "

$ws.Range("D19").Value = "42.323.70"
$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000104"
$ws.Range("E20").Value = "  -0.41%  "

$ws.Range("E21").Value = "  -2.44%  "

$ws.Range("E22").Value = "  -4.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.05"
$ws.Range("E23").Value = "  -9.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "228.66"
$ws.Range("E24").Value = "  -0.43%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.12"
$ws.Range("E25").Value = "  +5.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.91"
$ws.Range("E27").Value = "  -3.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.35"
$ws.Range("E28").Value = "  -7.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -2.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.19"
$ws.Range("E30").Value = "  -0.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.71"
$ws.Range("E31").Value = "  +3.12%  "

$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.31"
$ws.Range("E32").Value = "  +13.44%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.21"
$ws.Range("E33").Value = "  -1.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0793"
$ws.Range("E34").Value = "  -0.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.34"
$ws.Range("E35").Value = "  -4.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.122"
$ws.Range("E36").Value = "  -2.09%  "

$ws.Range("E37").Value = "  -2.15%  "

$ws.Range("E38").Value = "  +2.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0323"
$ws.Range("E39").Value = "  +6.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.63"
$ws.Range("E40").Value = "  -3.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.13"
$ws.Range("E41").Value = "  +0.24%  "

$ws.Range("E42").Value = "  -3.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "60.51"
$ws.Range("E43").Value = "  -6.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.197"
$ws.Range("E44").Value = "  -0.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.55"
$ws.Range("E45").Value = "  -1.63%  "

$ws.Range("E46").Value = "  -2.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.48"
$ws.Range("E47").Value = "  -3.83%  "

$ws.Range("E48").Value = "  -2.68%  "

$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("E50").Value = "  -2.62%  "

$ws.Range("E51").Value = "  +13.76%  "
